$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = "In Translation"
$wsOverview.Range("C2").Value = "In Translation"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("B2").Value = "In Translation"
$wsZhCn.Range("D3").Value = "2016-01-26 08:10:49"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("B2").Value = "In Translation"
$wsDeDe.Range("D3").Value = "2016-01-26 08:11:00"
